$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 3486.75
$ws.Range("J17").Value2 = 3486.75
$ws.Range("L17").Value2 = 10460.25
$ws.Range("N17").Value2 = -10796.25
$ws.Range("H19").Value2 = 7099.6
$ws.Range("H97").Value2 = 3611
$ws.Range("J97").Value2 = 3611
$ws.Range("L97").Value2 = 10833
$ws.Range("N97").Value2 = -11825
$ws.Range("H107").Value2 = 249.81818
$ws.Range("I107").Value2 = 250.33333
$ws.Range("J107").Value2 = 249.2
$ws.Range("K107").Value2 = 250.33333
$ws.Range("L107").Value2 = 249.2
$ws.Range("M107").Value2 = 1669.66667
$ws.Range("N107").Value2 = -4089.2
$ws.Range("H132").Value2 = 2579.8333
$ws.Range("I132").Value2 = 2122.25
$ws.Range("K132").Value2 = 6366.75
$ws.Range("M132").Value2 = -3836.75
$ws.Range("H137").Value2 = 1718.625
$ws.Range("I137").Value2 = 1722.9231
$ws.Range("K137").Value2 = 5168.7693
$ws.Range("M137").Value2 = -2618.7693
$ws.Range("H138").Value2 = 2831.3635
$ws.Range("I138").Value2 = 1722
$ws.Range("J138").Value2 = 3755.8333
$ws.Range("K138").Value2 = 5166
$ws.Range("L138").Value2 = 11267.4999
$ws.Range("M138").Value2 = -26
$ws.Range("N138").Value2 = -21547.4999
$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value2 = 0

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2003.8158
$ws.Range("I32").Value2 = 1962.3889
$ws.Range("K32").Value2 = 1962.3889
$ws.Range("M32").Value2 = -1675.3889
$ws.Range("H122").Value2 = 3725.6667
$ws.Range("I122").Value2 = 3538.5
$ws.Range("J122").Value2 = 4100
$ws.Range("K122").Value2 = 10615.5
$ws.Range("L122").Value2 = 12300
$ws.Range("M122").Value2 = -8165.5
$ws.Range("N122").Value2 = -17200
$ws.Range("H132").Value2 = 2644
$ws.Range("I132").Value2 = 1466.1666
$ws.Range("K132").Value2 = 4398.4998
$ws.Range("M132").Value2 = -1868.4998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 3057
$ws.Range("I86").Value2 = 2089.5
$ws.Range("J86").Value2 = 4024.5
$ws.Range("K86").Value2 = 2089.5
$ws.Range("L86").Value2 = 4024.5
$ws.Range("M86").Value2 = -966.5
$ws.Range("N86").Value2 = -6270.5
$ws.Range("H89").Value2 = 3057
$ws.Range("I89").Value2 = 2089.5
$ws.Range("J89").Value2 = 4024.5
$ws.Range("K89").Value2 = 10447.5
$ws.Range("L89").Value2 = 20122.5
$ws.Range("M89").Value2 = -4831.5
$ws.Range("N89").Value2 = -31354.5
$ws.Range("H123").Value2 = 0
$ws.Range("J123").Value2 = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value2 = 0
$ws.Range("H134").Value2 = 4353.2383
$ws.Range("J134").Value2 = 2967.7778
$ws.Range("L134").Value2 = 8903.3334
$ws.Range("N134").Value2 = -13973.3334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3753.1333
$ws.Range("I31").Value2 = 1260
$ws.Range("J31").Value2 = 4999.7
$ws.Range("K31").Value2 = 1260
$ws.Range("L31").Value2 = 4999.7
$ws.Range("M31").Value2 = -965
$ws.Range("N31").Value2 = -5589.7
$ws.Range("H34").Value2 = 3753.1333
$ws.Range("I34").Value2 = 1260
$ws.Range("J34").Value2 = 4999.7
$ws.Range("K34").Value2 = 1260
$ws.Range("L34").Value2 = 4999.7
$ws.Range("M34").Value2 = -1058
$ws.Range("N34").Value2 = -5403.7
$ws.Range("H58").Value2 = 1497.3334
$ws.Range("I58").Value2 = 1420.2727
$ws.Range("J58").Value2 = 2345
$ws.Range("K58").Value2 = 1420.2727
$ws.Range("L58").Value2 = 2345
$ws.Range("M58").Value2 = -1217.2727
$ws.Range("N58").Value2 = -2751
$ws.Range("H62").Value2 = 4503
$ws.Range("J62").Value2 = 4503
$ws.Range("L62").Value2 = 4503
$ws.Range("N62").Value2 = -5751
$ws.Range("H65").Value2 = 4503
$ws.Range("J65").Value2 = 4503
$ws.Range("L65").Value2 = 22515
$ws.Range("N65").Value2 = -28755
$ws.Range("H105").Value2 = 3750
$ws.Range("I105").Value2 = 3437.5
$ws.Range("K105").Value2 = 3437.5
$ws.Range("M105").Value2 = -1690.5
$ws.Range("H122").Value2 = 2546.6
$ws.Range("I122").Value2 = 1980.5
$ws.Range("J122").Value2 = 2924
$ws.Range("K122").Value2 = 5941.5
$ws.Range("L122").Value2 = 8772
$ws.Range("M122").Value2 = -3491.5
$ws.Range("N122").Value2 = -13672
$ws.Range("H132").Value2 = 2214.1052
$ws.Range("I132").Value2 = 1879.375
$ws.Range("K132").Value2 = 5638.125
$ws.Range("M132").Value2 = -3108.125
$ws.Range("H134").Value2 = 2424
$ws.Range("I134").Value2 = 2424
$ws.Range("K134").Value2 = 7272
$ws.Range("M134").Value2 = -4737
$ws.Range("H136").Value2 = 1497.3334
$ws.Range("I136").Value2 = 1420.2727
$ws.Range("J136").Value2 = 2345
$ws.Range("K136").Value2 = 4260.8181
$ws.Range("L136").Value2 = 7035
$ws.Range("M136").Value2 = -1710.8181
$ws.Range("N136").Value2 = -12135

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value2 = 2729.625
$ws.Range("I64").Value2 = 1034.5
$ws.Range("J64").Value2 = 4424.75
$ws.Range("K64").Value2 = 3103.5
$ws.Range("L64").Value2 = 13274.25
$ws.Range("M64").Value2 = -2833.5
$ws.Range("N64").Value2 = -13814.25
$ws.Range("H67").Value2 = 2729.625
$ws.Range("I67").Value2 = 1034.5
$ws.Range("J67").Value2 = 4424.75
$ws.Range("K67").Value2 = 3103.5
$ws.Range("L67").Value2 = 13274.25
$ws.Range("M67").Value2 = -2167.5
$ws.Range("N67").Value2 = -15146.25
$ws.Range("H121").Value2 = 638.2857
$ws.Range("I121").Value2 = 349
$ws.Range("J121").Value2 = 754
$ws.Range("K121").Value2 = 1047
$ws.Range("L121").Value2 = 2262
$ws.Range("M121").Value2 = 263
$ws.Range("N121").Value2 = -4882

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value2 = 1868.0714
$ws.Range("I107").Value2 = 1263
$ws.Range("J107").Value2 = 2674.8333
$ws.Range("K107").Value2 = 1263
$ws.Range("L107").Value2 = 2674.8333
$ws.Range("M107").Value2 = 657
$ws.Range("N107").Value2 = -6514.8333
$ws.Range("H126").Value2 = 3000
$ws.Range("I126").Value2 = 3000
$ws.Range("J126").Value2 = 0
$ws.Range("K126").Value2 = 9000
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value2 = 0
$ws.Range("M126").Value2 = -6530
$ws.Range("H132").Value2 = 3008.0833
$ws.Range("I132").Value2 = 2344.4443
$ws.Range("K132").Value2 = 7033.3329
$ws.Range("M132").Value2 = -4503.3329

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2750
$ws.Range("J7").Value2 = 3000
$ws.Range("L7").Value2 = 3000
$ws.Range("N7").Value2 = -3224
$ws.Range("H40").Value2 = 3495
$ws.Range("I40").Value2 = 3495
$ws.Range("K40").Value2 = 3495
$ws.Range("M40").Value2 = -3359
$ws.Range("H122").Value2 = 5741.25
$ws.Range("I122").Value2 = 4048.75
$ws.Range("K122").Value2 = 12146.25
$ws.Range("M122").Value2 = -9696.25
$ws.Range("H126").Value2 = 2750
$ws.Range("J126").Value2 = 3000
$ws.Range("L126").Value2 = 9000
$ws.Range("N126").Value2 = -13940
$ws.Range("H131").Value2 = 27081.5
$ws.Range("I131").Value2 = 29000
$ws.Range("K131").Value2 = 29000
$ws.Range("M131").Value2 = -23960
$ws.Range("H132").Value2 = 2427.25
$ws.Range("I132").Value2 = 1193
$ws.Range("K132").Value2 = 3579
$ws.Range("M132").Value2 = -1049
$ws.Range("H136").Value2 = 5871.636
$ws.Range("I136").Value2 = 3655.8572
$ws.Range("J136").Value2 = 9749.25
$ws.Range("K136").Value2 = 10967.5716
$ws.Range("L136").Value2 = 29247.75
$ws.Range("M136").Value2 = -8417.571599999999
$ws.Range("N136").Value2 = -34347.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 2100
$ws.Range("I122").Value2 = 2250
$ws.Range("J122").Value2 = 1500
$ws.Range("K122").Value2 = 6750
$ws.Range("L122").Value2 = 4500
$ws.Range("M122").Value2 = -4300
$ws.Range("N122").Value2 = -9400
$ws.Range("I126").Value2 = 1937.5
$ws.Range("J126").Value2 = 1915.6666
$ws.Range("K126").Value2 = 5812.5
$ws.Range("L126").Value2 = 5746.9998
$ws.Range("M126").Value2 = -3342.5
$ws.Range("N126").Value2 = -10686.9998
